$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 10-11, pushing the existing
# Mordialloc (x2) and Southbank rows down to rows 12-14.
$ws.Rows("10:11").Insert()

# Fill the newly inserted row 10 with the Melbourne venue entry.
$ws.Range("A10").Value = "Melbourne"
$ws.Range("B10").Value = "Melbourne Central Lion Hotel, 211 La Trobe Street"
$ws.Range("C10").Value = "28/12/2020 10:30pm-12.00am"
$ws.Range("D10").Value = "Case attended venue"

# Fill the newly inserted row 11 with the Moorabbin venue entry.
$ws.Range("A11").Value = "Moorabbin"
$ws.Range("B11").Value = "Grape and Grain Liquor Cellars, 14/16 Station St"
$ws.Range("C11").Value = "21/12/20 2pm - 10pm  22/12/20 10am - 6pm  24/12/20 1pm - 10pm  28/12/20 8.05pm - 8.47pm  29/12/20 12pm - 4pm"
$ws.Range("D11").Value = "Case's workplace"
